# Weekly refresh of the Frutilla (strawberry) price table.
# The whole data block (rows 2-30) shifts up by one record: a brand-new
# latest-week record is written at row 2, and the oldest week (previously
# rows 2-4) is appended as new rows 31-33 at the bottom, keeping the
# rolling 4-week x Calidad layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: convert an Excel date-serial (1900 system) string to a real date
# (no time-of-day component, matching the source cells which are midnight).
function SerialToDate([int]$serial) {
    $epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
    return $epoch.AddDays($serial)
}

# Final state for columns D, L, M, N, O, P, S for rows 2..33.
# (All other columns - A,B,C,E,F,G,H,I,J,K,Q,R,T - are identical on every
# data row, so they do not need to move.)
$data = @(
    @{ Row=2;  D=44708; L="Primera";  M=50;  N=6000; O=7000; P=6500; S=2167 },
    @{ Row=3;  D=44708; L="Segunda";  M=60;  N=4000; O=5000; P=4500; S=1500 },
    @{ Row=4;  D=44708; L="Tercera";  M=50;  N=3000; O=4000; P=3500; S=1167 },
    @{ Row=5;  D=44596; L="Especial"; M=100; N=8000; O=9000; P=8500; S=2833 },
    @{ Row=6;  D=44596; L="Primera";  M=130; N=6000; O=7000; P=6500; S=2167 },
    @{ Row=7;  D=44596; L="Segunda";  M=160; N=5000; O=6000; P=5750; S=1833 },
    @{ Row=8;  D=44596; L="Tercera";  M=100; N=4000; O=5000; P=4500; S=1500 },
    @{ Row=9;  D=44172; L="Especial"; M=100; N=6500; O=7000; P=6750; S=2250 },
    @{ Row=10; D=44172; L="Primera";  M=160; N=5500; O=6000; P=5750; S=1917 },
    @{ Row=11; D=44172; L="Segunda";  M=160; N=5000; O=5500; P=5250; S=1750 },
    @{ Row=12; D=44172; L="Tercera";  M=140; N=3500; O=4000; P=3750; S=1250 },
    @{ Row=13; D=44200; L="Especial"; M=50;  N=4500; O=5000; P=4750; S=1583 },
    @{ Row=14; D=44200; L="Primera";  M=80;  N=3500; O=4000; P=3750; S=1250 },
    @{ Row=15; D=44200; L="Segunda";  M=120; N=2500; O=3000; P=2750; S=917  },
    @{ Row=16; D=44322; L="Especial"; M=200; N=7000; O=7500; P=7250; S=2417 },
    @{ Row=17; D=44322; L="Primera";  M=160; N=6000; O=6500; P=6250; S=2083 },
    @{ Row=18; D=44322; L="Segunda";  M=100; N=5000; O=5500; P=5250; S=1750 },
    @{ Row=19; D=44249; L="Especial"; M=200; N=6000; O=7000; P=6500; S=2167 },
    @{ Row=20; D=44249; L="Primera";  M=160; N=4500; O=5000; P=4750; S=1583 },
    @{ Row=21; D=44334; L="Especial"; M=100; N=7000; O=8000; P=7500; S=2500 },
    @{ Row=22; D=44334; L="Primera";  M=160; N=6000; O=7000; P=6500; S=2167 },
    @{ Row=23; D=44334; L="Segunda";  M=120; N=6000; O=7000; P=6500; S=2167 },
    @{ Row=24; D=44334; L="Tercera";  M=70;  N=3500; O=4000; P=3750; S=1250 },
    @{ Row=25; D=44389; L="Especial"; M=100; N=7500; O=8000; P=7750; S=2583 },
    @{ Row=26; D=44389; L="Primera";  M=160; N=6000; O=7000; P=6500; S=2167 },
    @{ Row=27; D=44389; L="Segunda";  M=200; N=5500; O=6000; P=5750; S=1917 },
    @{ Row=28; D=44242; L="Especial"; M=50;  N=7000; O=8000; P=7500; S=2500 },
    @{ Row=29; D=44242; L="Primera";  M=90;  N=6000; O=7000; P=6500; S=2167 },
    @{ Row=30; D=44242; L="Segunda";  M=100; N=4000; O=5000; P=4500; S=1500 },
    @{ Row=31; D=44351; L="Especial"; M=160; N=7500; O=8000; P=7750; S=2583 },
    @{ Row=32; D=44351; L="Primera";  M=100; N=6000; O=6500; P=6250; S=2083 },
    @{ Row=33; D=44351; L="Segunda";  M=200; N=4500; O=5000; P=4750; S=1583 }
)

# Template row (row 2, before edits) holds the values common to every data
# row - copy it down to the three brand-new rows (31-33) first so all
# columns (A, B, C, E, F, G, H, I, J, K, Q, R, T) are populated, then
# overwrite the per-row fields (D, L, M, N, O, P, S) for every row.
$templateRow = 2
foreach ($newRow in 31..33) {
    $src = $ws.Range($ws.Cells.Item($templateRow, 1), $ws.Cells.Item($templateRow, 20))
    $dst = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 20))
    $src.Copy($dst)
}

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = (SerialToDate $entry.D)
    $ws.Cells.Item($r, 12).Value = $entry.L
    $ws.Cells.Item($r, 13).Value = $entry.M
    $ws.Cells.Item($r, 14).Value = $entry.N
    $ws.Cells.Item($r, 15).Value = $entry.O
    $ws.Cells.Item($r, 16).Value = $entry.P
    $ws.Cells.Item($r, 19).Value = $entry.S
}
